$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the stray "_GoBack" bookmark that currently sits at the end of
#    the "Problemes non resolus :" paragraph. (It will be re-created further
#    up the document as part of change #3.)
# ---------------------------------------------------------------------------
try {
    $oldGoBack = $d.Bookmarks.Item("_GoBack")
    $oldGoBack.Delete()
} catch {
    # no-op if it somehow is not present
}

# ---------------------------------------------------------------------------
# 2) "Auteur :" line -- merge " " + "Pelissier" + " Thomas" (which were split
#    across three runs / wrapped in a spellStart-spellEnd proofErr pair) into
#    a single run " Pelissier Thomas".
# ---------------------------------------------------------------------------
$rAuteur = $d.Content
$rAuteur.Find.Execute(
    " Pelissier Thomas", $false, $false, $false, $false, $false, $true,
    1, $false, " Pelissier Thomas", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) "Nom :" line -- change "... (PACKAGE : )" into "... (PACKAGE : Gestion)"
#    with "Gestion" as its own run and a "_GoBack" bookmark placed between
#    "Gestion" and the closing parenthesis. The whole paragraph is rebuilt
#    (deleted and re-created) so that the leftover gramStart/gramEnd proofing
#    marks that used to wrap ": )" are fully discarded rather than merely
#    shifted elsewhere.
# ---------------------------------------------------------------------------
$nomParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $txt = $d.Paragraphs.Item($i).Range.Text
    if ($txt -like "Nom*identifier (PACKAGE*") {
        $nomParaIndex = $i
        break
    }
}

if ($nomParaIndex -ne -1) {
    $nomPara = $d.Paragraphs.Item($nomParaIndex)

    # Delete the paragraph's content *and* its own trailing paragraph mark as
    # a single operation -- this is what makes the end-of-paragraph proofErr
    # markers disappear instead of migrating onto the following paragraph.
    $nomPara.Range.Delete()

    # Re-create an (empty) paragraph in the same spot, inheriting formatting
    # from its neighbours, then fill it back in.
    $afterPara = $d.Paragraphs.Item($nomParaIndex)
    $insPos = $afterPara.Range.Start
    $insRng = $d.Range($insPos, $insPos)
    $insRng.InsertParagraphBefore()

    $newPara = $d.Paragraphs.Item($nomParaIndex)
    $base = $newPara.Range.Start

    $fullText = "Nom" + [char]0x00A0 + ": S" + [char]0x2019 + "identifier (PACKAGE" + [char]0x00A0 + ": Gestion)"
    $fillRng = $d.Range($base, $base)
    $fillRng.InsertAfter($fullText)

    # "Nom :" (5 characters incl. the nbsp) keeps its original underline.
    $underlineRng = $d.Range($base, $base + 5)
    $underlineRng.Font.Underline = 1

    # Split "Gestion" away from the " ... (PACKAGE : " prefix -- use a
    # throwaway bookmark purely to force the run boundary, then remove it.
    $splitPos = $base + 30
    $splitRng = $d.Range($splitPos, $splitPos)
    $d.Bookmarks.Add("zzTmpSplit", $splitRng) | Out-Null
    $d.Bookmarks.Item("zzTmpSplit").Delete()

    # Place the real "_GoBack" bookmark right before the closing ")".
    $bmPos = $base + 37
    $bmRng = $d.Range($bmPos, $bmPos)
    $d.Bookmarks.Add("_GoBack", $bmRng) | Out-Null
}

# ---------------------------------------------------------------------------
# 4) Table cell -- merge "Demande le login" + " et le " into a single run
#    "Demande le login et le ".
# ---------------------------------------------------------------------------
$rLogin = $d.Content
$rLogin.Find.Execute(
    "Demande le login et le ", $false, $false, $false, $false, $false, $true,
    1, $false, "Demande le login et le ", 2) | Out-Null
